# Daily attendance processing - 2025-11-12 17:47:07
# Swap the order of the two comma-separated "Recorded By" values in column G
# (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("dnasr281@gmail.com, ")) {
        $parts = $val -split ", "
        if ($parts.Count -eq 2) {
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newVal
        }
    }
}
